# Generate Report for Handback
#
# Once a locale's localized content has round-tripped back in sync with
# en-US, this fills in the "handed back" columns for that locale:
#   - Status                     -> "Handed back: in sync with en-US"
#   - Latest Target File (I)     -> the source markdown doc (linked, like A)
#   - Latest Handback File (J)   -> the localized xliff that came back
#   - Latest Handback DateTime (K) -> when the handback happened
# and widens a few columns so the new / longer values stay readable.

$wb = $excel.ActiveWorkbook

$mdFileName  = "ee778694-5ccc-4a6d-9240-5f72d64a0014.md"
$mdHyperlink = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e4bf2a91e38440ef90106b4301830094ae54586/e2e/ee778694-5ccc-4a6d-9240-5f72d64a0014.md"
$newStatus   = "Handed back: in sync with en-US"

# ColumnWidth values land on Excel's internal 1/6-character pixel grid as
# (ColumnWidth + 5/6), rounded to the nearest 1/6. These inputs are chosen
# so the saved <col width="..."> comes out as close as possible to the
# target widths (30 is the nearest grid point to 29.9777047293527; 40
# lands on the grid exactly):
$wideColWidth   = 29.16666666666667   # -> saved width 30
$extraWideWidth = 39.16666666666667   # -> saved width 40

# ---------------------------------------------------------------------
# Overview sheet - the zh-cn / de-de roll-up cells show the very same
# Status text as the detail sheets, so they flip to the new text too.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E:E").ColumnWidth = $wideColWidth
$wsOverview.Range("F:F").ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdHyperlink, [System.Type]::Missing, [System.Type]::Missing, $mdFileName)
$wsZhCn.Range("J2").Value = "ee778694-5ccc-4a6d-9240-5f72d64a0014.56a6c5df332dc760f00010a0507bc54ab8f66aa3.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-28 20:57:41"
$wsZhCn.Range("C:C").ColumnWidth = $wideColWidth
$wsZhCn.Range("I:I").ColumnWidth = $extraWideWidth
$wsZhCn.Range("J:J").ColumnWidth = $extraWideWidth

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdHyperlink, [System.Type]::Missing, [System.Type]::Missing, $mdFileName)
$wsDeDe.Range("J2").Value = "ee778694-5ccc-4a6d-9240-5f72d64a0014.56a6c5df332dc760f00010a0507bc54ab8f66aa3.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-28 20:57:48"
$wsDeDe.Range("C:C").ColumnWidth = $wideColWidth
$wsDeDe.Range("I:I").ColumnWidth = $extraWideWidth
$wsDeDe.Range("J:J").ColumnWidth = $extraWideWidth

Write-Host "Handback report generated for zh-cn and de-de."
